# Updated fitting parameters and ready tu run detail tests.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Parameters")

# Update the fitted h_p_star value in K2
$ws.Range("K2").Value = 0.28537

# Leave the active selection on K3, as it was left after the edit
$ws.Activate()
$ws.Range("K3").Select()
